$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely
$ws.Rows.Item(11).Delete()

# ---- Elements sheet ----
$ws2 = $wb.Worksheets.Item("Elements")

# Row 2 (the root Extension element) gets a real Short/Definition instead of the defaults
$ws2.Cells.Item(2, 11).Value = "Medicare Coverage"
$ws2.Cells.Item(2, 12).Value = "Standard code for the type of Medicare coverage, if any, for the person"
